$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 "Modelo" - copy formatting from the existing header (E1)
# so it picks up the same bold/border/alignment style as the rest of row 1.
$ws.Range("F1").Value = "Modelo"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Updated metric values in row 2
$ws.Range("B2").Value = 0.1494536832641422
$ws.Range("C2").Value = 0.9984596127540851
$ws.Range("D2").Value = 0.3260639631584628

# New model-name cell F2
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor())])"
